$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 699659.0600000001
$ws.Range("J17").Value = 724459.4
$ws.Range("L17").Value = 2173378.2
$ws.Range("N17").Value = -2173714.2
$ws.Range("H19").Value = 1719.5555
$ws.Range("I19").Value = 1719.5555
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 1719.5555
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -1544.5555
$ws.Range("N19").ClearContents()
$ws.Range("H41").Value = 3185.182
$ws.Range("I41").Value = 1626.7142
$ws.Range("J41").Value = 5912.5
$ws.Range("K41").Value = 1626.7142
$ws.Range("L41").Value = 5912.5
$ws.Range("M41").Value = -1186.7142
$ws.Range("N41").Value = -6792.5
$ws.Range("H44").Value = 43993.332
$ws.Range("J44").Value = 43993.332
$ws.Range("L44").Value = 43993.332
$ws.Range("N44").Value = -44917.332
$ws.Range("H69").Value = 7297.409
$ws.Range("I69").Value = 5013
$ws.Range("J69").Value = 7406.1904
$ws.Range("K69").Value = 15039
$ws.Range("L69").Value = 22218.5712
$ws.Range("M69").Value = -14165
$ws.Range("N69").Value = -23966.5712
$ws.Range("H72").Value = 7297.409
$ws.Range("I72").Value = 5013
$ws.Range("J72").Value = 7406.1904
$ws.Range("K72").Value = 45117
$ws.Range("L72").Value = 66655.7136
$ws.Range("M72").Value = -40749
$ws.Range("N72").Value = -75391.7136
$ws.Range("H92").Value = 19232066
$ws.Range("I92").Value = 23810824
$ws.Range("J92").Value = 1283.2
$ws.Range("K92").Value = 23810824
$ws.Range("L92").Value = 1283.2
$ws.Range("M92").Value = -23809576
$ws.Range("N92").Value = -3779.2
$ws.Range("H96").Value = 1678.6666
$ws.Range("I96").Value = 1480
$ws.Range("K96").Value = 4440
$ws.Range("M96").Value = -3067
$ws.Range("H101").Value = 960
$ws.Range("I101").Value = 882.5
$ws.Range("J101").Value = 1089.1666
$ws.Range("K101").Value = 2647.5
$ws.Range("L101").Value = 3267.4998
$ws.Range("M101").Value = -1025.5
$ws.Range("N101").Value = -6511.4998
$ws.Range("H104").Value = 1499.5714
$ws.Range("I104").Value = 1307.8
$ws.Range("J104").Value = 1979
$ws.Range("K104").Value = 3923.4
$ws.Range("L104").Value = 5937
$ws.Range("M104").Value = -2176.4
$ws.Range("N104").Value = -9431
$ws.Range("H132").Value = 1337.72
$ws.Range("I132").Value = 1366.375
$ws.Range("J132").Value = 650
$ws.Range("K132").Value = 4099.125
$ws.Range("L132").Value = 1950
$ws.Range("M132").Value = -1569.125
$ws.Range("N132").Value = -7010
$ws.Range("H137").Value = 29266.9
$ws.Range("I137").Value = 30224.69
$ws.Range("K137").Value = 90674.06999999999
$ws.Range("M137").Value = -88124.06999999999
$ws.Range("H138").Value = 4035.77
$ws.Range("I138").Value = 3447.6
$ws.Range("J138").Value = 4101.122
$ws.Range("K138").Value = 10342.8
$ws.Range("L138").Value = 12303.366
$ws.Range("M138").Value = -5202.799999999999
$ws.Range("N138").Value = -22583.366

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2615.3242
$ws.Range("I61").Value = 1921.2727
$ws.Range("J61").Value = 3633.2666
$ws.Range("K61").Value = 1921.2727
$ws.Range("L61").Value = 3633.2666
$ws.Range("M61").Value = -1709.2727
$ws.Range("N61").Value = -4057.2666
$ws.Range("H63").Value = 5571.4287
$ws.Range("J63").Value = 5500
$ws.Range("L63").Value = 5500
$ws.Range("N63").Value = -6872
$ws.Range("H66").Value = 5571.4287
$ws.Range("J66").Value = 5500
$ws.Range("L66").Value = 27500
$ws.Range("N66").Value = -34364
$ws.Range("H132").Value = 2398.8367
$ws.Range("I132").Value = 2000.7561
$ws.Range("K132").Value = 6002.2683
$ws.Range("M132").Value = -3472.2683
$ws.Range("H136").Value = 2615.3242
$ws.Range("I136").Value = 1921.2727
$ws.Range("J136").Value = 3633.2666
$ws.Range("K136").Value = 5763.8181
$ws.Range("L136").Value = 10899.7998
$ws.Range("M136").Value = -3213.8181
$ws.Range("N136").Value = -15999.7998

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1106.6154
$ws.Range("I94").Value = 895.5
$ws.Range("K94").Value = 895.5
$ws.Range("M94").Value = -444.5
$ws.Range("H134").Value = 2978.2083
$ws.Range("I134").Value = 1843.15
$ws.Range("K134").Value = 5529.450000000001
$ws.Range("M134").Value = -2994.450000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 40.64706
$ws.Range("I7").Value = 36.8
$ws.Range("J7").Value = 69.5
$ws.Range("K7").Value = 36.8
$ws.Range("L7").Value = 69.5
$ws.Range("M7").Value = 76.2
$ws.Range("N7").Value = -295.5
$ws.Range("H22").Value = 817.625
$ws.Range("I22").Value = 820.1429000000001
$ws.Range("K22").Value = 820.1429000000001
$ws.Range("M22").Value = -470.1429000000001
$ws.Range("H31").Value = 5875.1626
$ws.Range("I31").Value = 1621.6111
$ws.Range("J31").Value = 8937.719999999999
$ws.Range("K31").Value = 1621.6111
$ws.Range("L31").Value = 8937.719999999999
$ws.Range("M31").Value = -1326.6111
$ws.Range("N31").Value = -9527.719999999999
$ws.Range("H34").Value = 5875.1626
$ws.Range("I34").Value = 1621.6111
$ws.Range("J34").Value = 8937.719999999999
$ws.Range("K34").Value = 1621.6111
$ws.Range("L34").Value = 8937.719999999999
$ws.Range("M34").Value = -1419.6111
$ws.Range("N34").Value = -9341.719999999999
$ws.Range("H58").Value = 2976.7715
$ws.Range("I58").Value = 2946.6765
$ws.Range("K58").Value = 2946.6765
$ws.Range("M58").Value = -2743.6765
$ws.Range("H105").Value = 1838.0476
$ws.Range("I105").Value = 1099.9286
$ws.Range("K105").Value = 1099.9286
$ws.Range("M105").Value = 647.0714
$ws.Range("H132").Value = 2675.7715
$ws.Range("I132").Value = 2843.6206
$ws.Range("J132").Value = 1864.5
$ws.Range("K132").Value = 8530.861800000001
$ws.Range("L132").Value = 5593.5
$ws.Range("M132").Value = -6000.861800000001
$ws.Range("N132").Value = -10653.5
$ws.Range("H134").Value = 1831.7354
$ws.Range("I134").Value = 1509
$ws.Range("K134").Value = 4527
$ws.Range("M134").Value = -1992
$ws.Range("H136").Value = 2976.7715
$ws.Range("I136").Value = 2946.6765
$ws.Range("K136").Value = 8840.029500000001
$ws.Range("M136").Value = -6290.029500000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 250000370
$ws.Range("I32").Value = 500
$ws.Range("K32").Value = 1500
$ws.Range("M32").Value = -1217
$ws.Range("H107").Value = 725.8333
$ws.Range("J107").Value = 557.2308
$ws.Range("L107").Value = 1671.6924
$ws.Range("N107").Value = -5511.6924
$ws.Range("H121").Value = 4774648.5
$ws.Range("I121").Value = 707.8
$ws.Range("J121").Value = 6479627.5
$ws.Range("K121").Value = 2123.4
$ws.Range("L121").Value = 19438882.5
$ws.Range("M121").Value = -813.3999999999996
$ws.Range("N121").Value = -19441502.5
$ws.Range("H129").Value = 1972
$ws.Range("I129").Value = 671
$ws.Range("K129").Value = 2013
$ws.Range("M129").Value = 2987
$ws.Range("H137").Value = 3768.3684
$ws.Range("J137").Value = 3869.6667
$ws.Range("L137").Value = 11609.0001
$ws.Range("N137").Value = -21809.0001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 5918.353
$ws.Range("I2").Value = 30.583334
$ws.Range("K2").Value = 30.583334
$ws.Range("M2").Value = 82.41666599999999
$ws.Range("H120").Value = 65960.60000000001
$ws.Range("J120").Value = 65960.60000000001
$ws.Range("L120").Value = 65960.60000000001
$ws.Range("N120").Value = -75636.60000000001
$ws.Range("H126").Value = 3449.8333
$ws.Range("J126").Value = 3574.75
$ws.Range("L126").Value = 10724.25
$ws.Range("N126").Value = -15664.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3682.1904
$ws.Range("I22").Value = 3234.4546
$ws.Range("J22").Value = 4174.7
$ws.Range("K22").Value = 3234.4546
$ws.Range("L22").Value = 4174.7
$ws.Range("M22").Value = -2939.4546
$ws.Range("N22").Value = -4764.7
$ws.Range("H27").Value = 3682.1904
$ws.Range("I27").Value = 3234.4546
$ws.Range("J27").Value = 4174.7
$ws.Range("K27").Value = 3234.4546
$ws.Range("L27").Value = 4174.7
$ws.Range("M27").Value = -3127.4546
$ws.Range("N27").Value = -4388.7
$ws.Range("H55").Value = 750.4545000000001
$ws.Range("I55").Value = 403.8
$ws.Range("K55").Value = 403.8
$ws.Range("M55").Value = -230.8
$ws.Range("H132").Value = 4790.3193
$ws.Range("I132").Value = 4708.575
$ws.Range("K132").Value = 14125.725
$ws.Range("M132").Value = -11595.725
$ws.Range("H136").Value = 1703.6471
$ws.Range("I136").Value = 1290.8572
$ws.Range("K136").Value = 3872.5716
$ws.Range("M136").Value = -1322.5716

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4592.091
$ws.Range("I62").Value = 3102.1667
$ws.Range("K62").Value = 3102.1667
$ws.Range("M62").Value = -2478.1667
$ws.Range("H65").Value = 4592.091
$ws.Range("I65").Value = 3102.1667
$ws.Range("K65").Value = 15510.8335
$ws.Range("M65").Value = -12390.8335
$ws.Range("H81").Value = 4305.905
$ws.Range("I81").Value = 3221.6365
$ws.Range("J81").Value = 5498.6
$ws.Range("K81").Value = 6443.273
$ws.Range("L81").Value = 10997.2
$ws.Range("M81").Value = -5382.273
$ws.Range("N81").Value = -13119.2
$ws.Range("H84").Value = 4305.905
$ws.Range("I84").Value = 3221.6365
$ws.Range("J84").Value = 5498.6
$ws.Range("K84").Value = 32216.365
$ws.Range("L84").Value = 54986
$ws.Range("M84").Value = -26912.365
$ws.Range("N84").Value = -65594
$ws.Range("H98").Value = 57500
$ws.Range("J98").Value = 57500
$ws.Range("L98").Value = 57500
$ws.Range("N98").Value = -63490
$ws.Range("H136").Value = 30824.371
$ws.Range("I136").Value = 1637.1666
$ws.Range("K136").Value = 4911.4998
$ws.Range("M136").Value = -2361.4998
